$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 123, shifting the existing
# rows 123-156 down to 124-157 (the last existing row ends up at 157).
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new weekly record.
$ws.Cells.Item(123, 1).Value = 4
$ws.Cells.Item(123, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(123, 3).Value = "Los Lagos"
$ws.Cells.Item(123, 4).Value = 44642
$ws.Cells.Item(123, 5).Value = 10
$ws.Cells.Item(123, 6).Value = 100112009
$ws.Cells.Item(123, 7).Value = "Acelga"
$ws.Cells.Item(123, 8).Value = "Sin especificar"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 80
$ws.Cells.Item(123, 11).Value = 10000
$ws.Cells.Item(123, 12).Value = 10000
$ws.Cells.Item(123, 13).Value = 10000
$ws.Cells.Item(123, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(123, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(123, 16).Value = 833
$ws.Cells.Item(123, 17).Value = 12
$ws.Cells.Item(123, 18).Value = "Hortaliza"
